$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033187354849678
$ws.Cells.Item(2, 4).Value = 1.042153669889568
$ws.Cells.Item(2, 5).Value = 1.032510211636113
$ws.Cells.Item(2, 6).Value = 1.052534165196856
$ws.Cells.Item(2, 9).Value = 1.036090546296486
$ws.Cells.Item(2, 10).Value = 1.038313457740253
$ws.Cells.Item(2, 11).Value = 1.044931161085016
$ws.Cells.Item(2, 12).Value = 1.035315234379497
$ws.Cells.Item(2, 13).Value = 1.055282631970104
$ws.Cells.Item(2, 14).Value = 1.016657186976145

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034374012609437
$ws.Cells.Item(3, 4).Value = 1.043074659486301
$ws.Cells.Item(3, 5).Value = 1.033525806109813
$ws.Cells.Item(3, 6).Value = 1.053597133454186
$ws.Cells.Item(3, 9).Value = 1.036330149038658
$ws.Cells.Item(3, 10).Value = 1.039141521538604
$ws.Cells.Item(3, 11).Value = 1.045662635342709
$ws.Cells.Item(3, 12).Value = 1.036139062976453
$ws.Cells.Item(3, 13).Value = 1.056157823300814
$ws.Cells.Item(3, 14).Value = 1.016938930208505

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.035141810305556
$ws.Cells.Item(4, 4).Value = 1.043670364051999
$ws.Cells.Item(4, 5).Value = 1.034183242216883
$ws.Cells.Item(4, 6).Value = 1.054285006378844
$ws.Cells.Item(4, 9).Value = 1.03648374834845
$ws.Cells.Item(4, 10).Value = 1.039676782418179
$ws.Cells.Item(4, 11).Value = 1.046135084700695
$ws.Cells.Item(4, 12).Value = 1.036671810007778
$ws.Cells.Item(4, 13).Value = 1.056723592954004
$ws.Cells.Item(4, 14).Value = 1.017120866469662

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035464582118876
$ws.Cells.Item(5, 4).Value = 1.043920741694654
$ws.Cells.Item(5, 5).Value = 1.034459695536843
$ws.Cells.Item(5, 6).Value = 1.054574203223369
$ws.Cells.Item(5, 9).Value = 1.0365479768909
$ws.Cells.Item(5, 10).Value = 1.039901674596635
$ws.Cells.Item(5, 11).Value = 1.04633349605155
$ws.Cells.Item(5, 12).Value = 1.036895699516599
$ws.Cells.Item(5, 13).Value = 1.056961314137956
$ws.Cells.Item(5, 14).Value = 1.0171972638803

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035518776362455
$ws.Cells.Item(6, 4).Value = 1.043962777893308
$ws.Cells.Item(6, 5).Value = 1.034506117219457
$ws.Cells.Item(6, 6).Value = 1.054622761540422
$ws.Cells.Item(6, 9).Value = 1.036558740939493
$ws.Cells.Item(6, 10).Value = 1.039939427289442
$ws.Cells.Item(6, 11).Value = 1.046366798099277
$ws.Cells.Item(6, 12).Value = 1.036933287029397
$ws.Cells.Item(6, 13).Value = 1.057001221052431
$ws.Cells.Item(6, 14).Value = 1.017210086149229

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.03514612323713
$ws.Cells.Item(7, 4).Value = 1.043673709832237
$ws.Cells.Item(7, 5).Value = 1.034186935935089
$ws.Cells.Item(7, 6).Value = 1.05428887058157
$ws.Cells.Item(7, 9).Value = 1.03648460792738
$ws.Cells.Item(7, 10).Value = 1.039679787954414
$ws.Cells.Item(7, 11).Value = 1.046137736693307
$ws.Cells.Item(7, 12).Value = 1.036674801933636
$ws.Cells.Item(7, 13).Value = 1.05672676990012
$ws.Cells.Item(7, 14).Value = 1.017121887643599

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033588402457827
$ws.Cells.Item(8, 4).Value = 1.042464971919816
$ws.Cells.Item(8, 5).Value = 1.032853379158938
$ws.Cells.Item(8, 6).Value = 1.052893387810938
$ws.Cells.Item(8, 9).Value = 1.036171819295818
$ws.Cells.Item(8, 10).Value = 1.038593420386545
$ws.Cells.Item(8, 11).Value = 1.045178545213209
$ws.Cells.Item(8, 12).Value = 1.035593718739106
$ws.Cells.Item(8, 13).Value = 1.055578518508745
$ws.Cells.Item(8, 14).Value = 1.016752480181088

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.030843057435372
$ws.Cells.Item(9, 4).Value = 1.040333192083163
$ws.Cells.Item(9, 5).Value = 1.03050559280833
$ws.Cells.Item(9, 6).Value = 1.050434821734592
$ws.Cells.Item(9, 9).Value = 1.03560961668347
$ws.Cells.Item(9, 10).Value = 1.03667484877815
$ws.Cells.Item(9, 11).Value = 1.043481698923875
$ws.Cells.Item(9, 12).Value = 1.033686198842673
$ws.Cells.Item(9, 13).Value = 1.0535510192778
$ws.Cells.Item(9, 14).Value = 1.016098695436681

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.029012437001407
$ws.Cells.Item(10, 4).Value = 1.038910756076721
$ws.Cells.Item(10, 5).Value = 1.028941787080602
$ws.Cells.Item(10, 6).Value = 1.048796054841862
$ws.Cells.Item(10, 9).Value = 1.035227394338901
$ws.Cells.Item(10, 10).Value = 1.035392900934258
$ws.Cells.Item(10, 11).Value = 1.042345986997178
$ws.Cells.Item(10, 12).Value = 1.032412794827768
$ws.Cells.Item(10, 13).Value = 1.052196549079663
$ws.Cells.Item(10, 14).Value = 1.015660921204388

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.028219644033878
$ws.Cells.Item(11, 4).Value = 1.038294523920311
$ws.Cells.Item(11, 5).Value = 1.028264961996579
$ws.Cells.Item(11, 6).Value = 1.048086509346878
$ws.Cells.Item(11, 9).Value = 1.035060126161848
$ws.Cells.Item(11, 10).Value = 1.034837105719147
$ws.Cells.Item(11, 11).Value = 1.041853142084248
$ws.Cells.Item(11, 12).Value = 1.031860979495899
$ws.Cells.Item(11, 13).Value = 1.051609377013835
$ws.Cells.Item(11, 14).Value = 1.015470903328261

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027925145110434
$ws.Cells.Item(12, 4).Value = 1.038065581021737
$ws.Cells.Item(12, 5).Value = 1.028013605236704
$ws.Cells.Item(12, 6).Value = 1.04782295955339
$ws.Cells.Item(12, 9).Value = 1.034997730180629
$ws.Cells.Item(12, 10).Value = 1.034630552022045
$ws.Cells.Item(12, 11).Value = 1.041669915519407
$ws.Cells.Item(12, 12).Value = 1.031655946345918
$ws.Cells.Item(12, 13).Value = 1.051391172709417
$ws.Cells.Item(12, 14).Value = 1.015400253098482

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027988317079915
$ws.Cells.Item(13, 4).Value = 1.03811469217932
$ws.Cells.Item(13, 5).Value = 1.028067520037126
$ws.Cells.Item(13, 6).Value = 1.047879491573821
$ws.Cells.Item(13, 9).Value = 1.035011126331515
$ws.Cells.Item(13, 10).Value = 1.034674863330549
$ws.Cells.Item(13, 11).Value = 1.04170922558162
$ws.Cells.Item(13, 12).Value = 1.031699929578438
$ws.Cells.Item(13, 13).Value = 1.05143798292044
$ws.Cells.Item(13, 14).Value = 1.015415410944249

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.028195301086732
$ws.Cells.Item(14, 4).Value = 1.038275600387672
$ws.Cells.Item(14, 5).Value = 1.028244183830499
$ws.Cells.Item(14, 6).Value = 1.04806472410022
$ws.Cells.Item(14, 9).Value = 1.035054973900389
$ws.Cells.Item(14, 10).Value = 1.034820034106307
$ws.Cells.Item(14, 11).Value = 1.04183799982812
$ws.Cells.Item(14, 12).Value = 1.031844032709544
$ws.Cells.Item(14, 13).Value = 1.051591342271871
$ws.Cells.Item(14, 14).Value = 1.015465064769674

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.028322828014185
$ws.Cells.Item(15, 4).Value = 1.038374735021954
$ws.Cells.Item(15, 5).Value = 1.028353038322889
$ws.Cells.Item(15, 6).Value = 1.048178852878487
$ws.Cells.Item(15, 9).Value = 1.035081954695941
$ws.Cells.Item(15, 10).Value = 1.034909464457855
$ws.Cells.Item(15, 11).Value = 1.041917320403007
$ws.Cells.Item(15, 12).Value = 1.031932810851013
$ws.Cells.Item(15, 13).Value = 1.051685818433552
$ws.Cells.Item(15, 14).Value = 1.015495648963627

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.029065049317717
$ws.Cells.Item(16, 4).Value = 1.038951646838311
$ws.Cells.Item(16, 5).Value = 1.028986712258564
$ws.Cells.Item(16, 6).Value = 1.048843146060415
$ws.Cells.Item(16, 9).Value = 1.035238458207077
$ws.Cells.Item(16, 10).Value = 1.035429772333607
$ws.Cells.Item(16, 11).Value = 1.042378672805083
$ws.Cells.Item(16, 12).Value = 1.032449408012421
$ws.Cells.Item(16, 13).Value = 1.052235503424459
$ws.Cells.Item(16, 14).Value = 1.015673522382005

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.029530591029007
$ws.Cells.Item(17, 4).Value = 1.039313445525112
$ws.Cells.Item(17, 5).Value = 1.029384282422311
$ws.Cells.Item(17, 6).Value = 1.049259852870362
$ws.Cells.Item(17, 9).Value = 1.035336156473821
$ws.Cells.Item(17, 10).Value = 1.035755958752392
$ws.Cells.Item(17, 11).Value = 1.042667778945483
$ws.Cells.Item(17, 12).Value = 1.032773341941658
$ws.Cells.Item(17, 13).Value = 1.052580124434974
$ws.Cells.Item(17, 14).Value = 1.01578497471222

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029802121969525
$ws.Cells.Item(18, 4).Value = 1.03952444690078
$ws.Cells.Item(18, 5).Value = 1.029616208639474
$ws.Cells.Item(18, 6).Value = 1.049502915909051
$ws.Cells.Item(18, 9).Value = 1.035392972093103
$ws.Cells.Item(18, 10).Value = 1.035946149972919
$ws.Cells.Item(18, 11).Value = 1.042836306068735
$ws.Cells.Item(18, 12).Value = 1.03296224636457
$ws.Cells.Item(18, 13).Value = 1.052781070712591
$ws.Cells.Item(18, 14).Value = 1.015849938722708

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029894705084268
$ws.Cells.Item(19, 4).Value = 1.039596387900553
$ws.Cells.Item(19, 5).Value = 1.02969529464502
$ws.Cells.Item(19, 6).Value = 1.049585794995097
$ws.Cells.Item(19, 9).Value = 1.035412315900797
$ws.Cells.Item(19, 10).Value = 1.036010988792777
$ws.Cells.Item(19, 11).Value = 1.042893751910206
$ws.Cells.Item(19, 12).Value = 1.033026651002677
$ws.Cells.Item(19, 13).Value = 1.052849577139843
$ws.Cells.Item(19, 14).Value = 1.015872082272544

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.029480644029449
$ws.Cells.Item(20, 4).Value = 1.039274631023628
$ws.Cells.Item(20, 5).Value = 1.029341623781462
$ws.Cells.Item(20, 6).Value = 1.049215143660515
$ws.Cells.Item(20, 9).Value = 1.035325691972828
$ws.Cells.Item(20, 10).Value = 1.035720969056014
$ws.Cells.Item(20, 11).Value = 1.0426367713127
$ws.Cells.Item(20, 12).Value = 1.0327385911182
$ws.Cells.Item(20, 13).Value = 1.052543156623736
$ws.Cells.Item(20, 14).Value = 1.0157730215102

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.028134350032169
$ws.Cells.Item(21, 4).Value = 1.038228218228379
$ws.Cells.Item(21, 5).Value = 1.028192159477586
$ws.Cells.Item(21, 6).Value = 1.048010177550706
$ws.Cells.Item(21, 9).Value = 1.035042069204834
$ws.Cells.Item(21, 10).Value = 1.03477728788157
$ws.Cells.Item(21, 11).Value = 1.041800083499852
$ws.Cells.Item(21, 12).Value = 1.031801599712502
$ws.Cells.Item(21, 13).Value = 1.051546184598871
$ws.Cells.Item(21, 14).Value = 1.015450444865687

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.027287763051361
$ws.Cells.Item(22, 4).Value = 1.03757002569821
$ws.Cells.Item(22, 5).Value = 1.027469712602504
$ws.Cells.Item(22, 6).Value = 1.047252607647445
$ws.Cells.Item(22, 9).Value = 1.034862210236106
$ws.Cells.Item(22, 10).Value = 1.034183340808354
$ws.Cells.Item(22, 11).Value = 1.041273086945352
$ws.Cells.Item(22, 12).Value = 1.031212103217743
$ws.Cells.Item(22, 13).Value = 1.050918755611266
$ws.Cells.Item(22, 14).Value = 1.01524722805662

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.027736566713448
$ws.Cells.Item(23, 4).Value = 1.037918971916182
$ws.Cells.Item(23, 5).Value = 1.027852670334072
$ws.Cells.Item(23, 6).Value = 1.04765420622723
$ws.Cells.Item(23, 9).Value = 1.034957702359889
$ws.Cells.Item(23, 10).Value = 1.034498262213479
$ws.Cells.Item(23, 11).Value = 1.041552546889699
$ws.Cells.Item(23, 12).Value = 1.031524642074757
$ws.Cells.Item(23, 13).Value = 1.051251424031141
$ws.Cells.Item(23, 14).Value = 1.015354995124267

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.029503212964111
$ws.Cells.Item(24, 4).Value = 1.03929216971754
$ws.Cells.Item(24, 5).Value = 1.029360899291223
$ws.Cells.Item(24, 6).Value = 1.049235345812195
$ws.Cells.Item(24, 9).Value = 1.035330420955942
$ws.Cells.Item(24, 10).Value = 1.035736779602866
$ws.Cells.Item(24, 11).Value = 1.042650782647088
$ws.Cells.Item(24, 12).Value = 1.032754293645137
$ws.Cells.Item(24, 13).Value = 1.052559860987951
$ws.Cells.Item(24, 14).Value = 1.015778422784303

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.031552856904883
$ws.Cells.Item(25, 4).Value = 1.040884527040907
$ws.Cells.Item(25, 5).Value = 1.031112305191769
$ws.Cells.Item(25, 6).Value = 1.051070369129622
$ws.Cells.Item(25, 9).Value = 1.035756267076748
$ws.Cells.Item(25, 10).Value = 1.03717135371039
$ws.Cells.Item(25, 11).Value = 1.043921162371347
$ws.Cells.Item(25, 12).Value = 1.034179639811785
$ws.Cells.Item(25, 13).Value = 1.0540756696257365
$ws.Cells.Item(25, 14).Value = 1.016268051943777
